$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.936.85'
$ws.Range("E2").Value = '  +1.98%  '
$ws.Range("D3").Value = '3.734.33'
$ws.Range("E3").Value = '  -1.97%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '601.42'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.42%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.87'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -2.25%  '
$ws.Range("D7").Value = '3.731.78'
$ws.Range("E7").Value = '  -1.89%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.532'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +2.07%  '
$ws.Range("E10").Value = '  +3.97%  '
$ws.Range("E11").Value = '  +1.45%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.459'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.42%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.24'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +1.37%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000244'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.12%  '
$ws.Range("D15").Value = '4.365.20'
$ws.Range("E15").Value = '  -1.52%  '
$ws.Range("D16").Value = '3.743.71'
$ws.Range("E16").Value = '  -1.24%  '
$ws.Range("D17").Value = '68.882.75'
$ws.Range("E17").Value = '  +1.74%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.26'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.65%  '
$ws.Range("E19").Value = '  +0.48%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.23'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +7.44%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '497.81'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +2.13%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.68'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +4.89%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.723'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.82'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +1.09%  '
$ws.Range("E25").Value = '  +1.46%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.30'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.55%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.26'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.66%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.10'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.87%  '
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("E30").Value = '  +1.09%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.41'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.52%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.91'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +1.33%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.69'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -3.08%  '
$ws.Range("D34").Value = '3.882.03'
$ws.Range("E34").Value = '  -1.41%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.108'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.43%  '
$ws.Range("D36").Value = '3.669.43'
$ws.Range("E36").Value = '  -1.86%  '
$ws.Range("E37").Value = '  +0.33%  '
$ws.Range("E38").Value = '  +0.26%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.77'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.73%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.133'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -1.54%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.323'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.10%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '436.35'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -3.27%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '49.00'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +0.21%  '
$ws.Range("E44").Value = '  -0.63%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.88'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +1.62%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.38'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.59%  '
$ws.Range("E47").Value = '  +0.01%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '40.46'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.36%  '
$ws.Range("E49").Value = '  +1.77%  '
$ws.Range("E50").Value = '  +0.79%  '
$ws.Range("D51").Value = '2.744.30'
$ws.Range("E51").Value = '  -3.26%  '
